# Weekly update: add a new week of "Sandia" (Hortaliza) price data for the
# "Terminal La Palmera de La Serena" market. This inserts three new rows
# (Extra / Primera / Segunda) at the top of that market's block (row 120),
# pushing all the existing rows for this market down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at row 120 (one at a time so each subsequent
# insert pushes the previously-inserted blank rows down, ending up with three
# consecutive new rows at 120-122 and the old row 120 now at row 123).
$ws.Rows.Item(120).Insert()
$ws.Rows.Item(120).Insert()
$ws.Rows.Item(120).Insert()

# New row 120: Sandia, "Extra" quality, week of 2023-01-17 (serial 44943)
$ws.Cells.Item(120, 1).Value = 8
$ws.Cells.Item(120, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44943
$ws.Cells.Item(120, 5).Value = 4
$ws.Cells.Item(120, 6).Value = 100112028
$ws.Cells.Item(120, 7).Value = "Sandia"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Extra"
$ws.Cells.Item(120, 10).Value = 2000
$ws.Cells.Item(120, 11).Value = 3300
$ws.Cells.Item(120, 12).Value = 3500
$ws.Cells.Item(120, 13).Value = 3400
$ws.Cells.Item(120, 14).Value = "$/unidad"
$ws.Cells.Item(120, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(120, 16).Value = 3400
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# New row 121: Sandia, "Primera" quality, week of 2023-01-17
$ws.Cells.Item(121, 1).Value = 8
$ws.Cells.Item(121, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(121, 3).Value = "Coquimbo"
$ws.Cells.Item(121, 4).Value = 44943
$ws.Cells.Item(121, 5).Value = 4
$ws.Cells.Item(121, 6).Value = 100112028
$ws.Cells.Item(121, 7).Value = "Sandia"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 1800
$ws.Cells.Item(121, 11).Value = 2800
$ws.Cells.Item(121, 12).Value = 3000
$ws.Cells.Item(121, 13).Value = 2900
$ws.Cells.Item(121, 14).Value = "$/unidad"
$ws.Cells.Item(121, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(121, 16).Value = 2900
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = "Hortaliza"

# New row 122: Sandia, "Segunda" quality, week of 2023-01-17
$ws.Cells.Item(122, 1).Value = 8
$ws.Cells.Item(122, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(122, 3).Value = "Coquimbo"
$ws.Cells.Item(122, 4).Value = 44943
$ws.Cells.Item(122, 5).Value = 4
$ws.Cells.Item(122, 6).Value = 100112028
$ws.Cells.Item(122, 7).Value = "Sandia"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Segunda"
$ws.Cells.Item(122, 10).Value = 1400
$ws.Cells.Item(122, 11).Value = 2300
$ws.Cells.Item(122, 12).Value = 2500
$ws.Cells.Item(122, 13).Value = 2400
$ws.Cells.Item(122, 14).Value = "$/unidad"
$ws.Cells.Item(122, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(122, 16).Value = 2400
$ws.Cells.Item(122, 17).Value = 1
$ws.Cells.Item(122, 18).Value = "Hortaliza"
